$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $escaped = $val -replace '"', '""'
    $r = $ws.Range($cellAddr)
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

Set-TextValue "D2" "29.209.29"
Set-TextValue "E2" "  +0.40%  "
Set-TextValue "D3" "1.857.47"
Set-TextValue "E3" "  +0.79%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "0.6996"
Set-TextValue "E5" "  +1.36%  "
Set-TextValue "D6" "237.33"
Set-TextValue "E6" "  +0.30%  "
Set-TextValue "D8" "0.07730"
Set-TextValue "E8" "  +2.74%  "
Set-TextValue "D9" "0.3042"
Set-TextValue "E9" "  +0.42%  "
Set-TextValue "D10" "23.26"
Set-TextValue "E10" "  -0.03%  "
Set-TextValue "D11" "0.08168"
Set-TextValue "E11" "  +1.18%  "
Set-TextValue "D12" "1.838.55"
Set-TextValue "E12" "  +2.70%  "
Set-TextValue "D13" "0.7180"
Set-TextValue "E13" "  -0.24%  "
Set-TextValue "D14" "5.160"
Set-TextValue "E14" "  -0.31%  "
Set-TextValue "D15" "89.09"
Set-TextValue "E15" "  +0.58%  "
Set-TextValue "D16" "29.205.70"
Set-TextValue "E16" "  +0.40%  "
Set-TextValue "D17" "5.763"
Set-TextValue "E17" "  -0.12%  "
Set-TextValue "D18" "13.33"
Set-TextValue "E18" "  +2.89%  "
Set-TextValue "D19" "0.000007726"
Set-TextValue "E19" "  +0.96%  "
Set-TextValue "D20" "236.57"
Set-TextValue "E20" "  -1.69%  "
Set-TextValue "D21" "0.9995"
Set-TextValue "E21" "  -0.11%  "
Set-TextValue "D22" "2.109.38"
Set-TextValue "E22" "  +0.93%  "
Set-TextValue "D23" "1.001"
Set-TextValue "E23" "  -0.04%  "
Set-TextValue "D24" "7.425"
Set-TextValue "E24" "  -2.26%  "
Set-TextValue "B25" "Stellar"
Set-TextValue "C25" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D25" "0.1478"
Set-TextValue "E25" "  +1.44%  "
Set-TextValue "B26" "Monero"
Set-TextValue "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "162.47"
Set-TextValue "E26" "  +0.57%  "
Set-TextValue "D27" "9.010"
Set-TextValue "E27" "  +0.32%  "
Set-TextValue "D28" "18.01"
Set-TextValue "E28" "  +0.07%  "
Set-TextValue "D29" "2.059"
Set-TextValue "E29" "  +7.16%  "
Set-TextValue "D30" "1.424"
Set-TextValue "E30" "  +3.36%  "
Set-TextValue "D31" "4.436"
Set-TextValue "E31" "  +0.49%  "
Set-TextValue "D32" "1.484"
Set-TextValue "E32" "  -0.13%  "
Set-TextValue "D33" "4.026"
Set-TextValue "E33" "  -0.04%  "
Set-TextValue "D34" "0.05217"
Set-TextValue "E34" "  +0.60%  "
Set-TextValue "D35" "1.167"
Set-TextValue "E35" "  -0.94%  "
Set-TextValue "E36" "  -0.35%  "
Set-TextValue "D37" "0.9998"
Set-TextValue "E37" "  -0.05%  "
Set-TextValue "D38" "2.665"
Set-TextValue "E38" "  +0.17%  "
Set-TextValue "D39" "0.01847"
Set-TextValue "E39" "  -0.41%  "
Set-TextValue "D40" "2.725"
Set-TextValue "E40" "  +1.94%  "
Set-TextValue "D41" "0.9355"
Set-TextValue "E41" "  +2.24%  "
Set-TextValue "D42" "1.141.41"
Set-TextValue "E42" "  +8.71%  "
Set-TextValue "D43" "0.4274"
Set-TextValue "E43" "  +0.28%  "
Set-TextValue "D44" "5.889"
Set-TextValue "E44" "  -0.27%  "
Set-TextValue "D45" "70.61"
Set-TextValue "E45" "  +1.72%  "
Set-TextValue "D46" "1.001"
Set-TextValue "E46" "  +0.02%  "
Set-TextValue "D47" "103.13"
Set-TextValue "E47" "  +0.83%  "
Set-TextValue "D48" "1.792"
Set-TextValue "E48" "  +3.05%  "
Set-TextValue "D49" "2.005.50"
Set-TextValue "E49" "  +0.57%  "
Set-TextValue "D50" "9.159"
Set-TextValue "E50" "  -0.63%  "
Set-TextValue "D51" "6.967"
Set-TextValue "E51" "  -2.33%  "